# Add a new "2022" column (S) to the right of the existing "2021" column (R),
# mirroring the formatting of column R for every data row, and move the
# active-cell selection as recorded in the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 holds the year headers; rows 5-34 hold the per-region data values
# (rows 5,8,11,... are "Women", 6,9,12,... are "Men", matching the existing
# R-column pattern). "-" denotes the existing placeholder shared string used
# throughout the sheet for missing values.
$values = @(
    2022,               # row 4  (header year)
    0.5,                # row 5
    0.2,                # row 6
    0.7,                # row 7
    0.2,                # row 8
    "-",                # row 9
    0.4,                # row 10
    0.5,                # row 11
    0.3,                # row 12
    0.6,                # row 13
    0.7,                # row 14
    0.4,                # row 15
    1.1000000000000001, # row 16
    "-",                # row 17
    "-",                # row 18
    "-",                # row 19
    0.4,                # row 20
    0.4,                # row 21
    0.4,                # row 22
    0.4,                # row 23
    "-",                # row 24
    0.7,                # row 25
    1,                  # row 26
    0.4,                # row 27
    1.7,                # row 28
    0.3,                # row 29
    0,                  # row 30
    0.6,                # row 31
    "-",                # row 32
    "-",                # row 33
    "-"                 # row 34
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 4 + $i

    # Write the new value into column S ...
    $ws.Range("S$row").Value = $values[$i]

    # ... then clone column R's formatting (number format/font/borders/
    # alignment) onto the new cell so it matches the rest of the table.
    $ws.Range("R$row").Copy()
    $ws.Range("S$row").PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# The recorded selection after the edit.
$ws.Range("T6").Select()
